# W14 Wednesday Commit 1
# updated chapter 6 report 1; garden excel sheet; garden observation sheet; timecard.
#
# Fills in the 6/26 and 6/27 (r17, r18) observation rows that were still
# blank on all four "Garden Practicum" sheets, adds a new "Leaf 7" column
# to the "Cardoon (2)" sheet, and leaves each sheet's selection where the
# user last clicked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Pansies Alive
# ---------------------------------------------------------------------
$wsPansiesAlive = $wb.Worksheets.Item("Pansies Alive")

$wsPansiesAlive.Range("B17").Value = 2
$wsPansiesAlive.Range("C17").Value = 5
$wsPansiesAlive.Range("D17").Value = 3
$wsPansiesAlive.Range("E17").Value = 1
$wsPansiesAlive.Range("F17").Value = 3

$wsPansiesAlive.Range("B18").Value = 2
$wsPansiesAlive.Range("C18").Value = 5
$wsPansiesAlive.Range("D18").Value = 3
$wsPansiesAlive.Range("E18").Value = 1
$wsPansiesAlive.Range("F18").Value = 3

$wsPansiesAlive.Range("F18").Select() | Out-Null

# ---------------------------------------------------------------------
# Pansies Dead
# ---------------------------------------------------------------------
$wsPansiesDead = $wb.Worksheets.Item("Pansies Dead")

$wsPansiesDead.Range("B17").Value = 0
$wsPansiesDead.Range("C17").Value = 0
$wsPansiesDead.Range("D17").Value = 0
$wsPansiesDead.Range("E17").Value = 0
$wsPansiesDead.Range("F17").Value = 0

$wsPansiesDead.Range("B18").Value = 0
$wsPansiesDead.Range("C18").Value = 0
$wsPansiesDead.Range("D18").Value = 0
$wsPansiesDead.Range("E18").Value = 0
$wsPansiesDead.Range("F18").Value = 0

$wsPansiesDead.Activate()
$wsPansiesDead.Range("F23").Select() | Out-Null

# ---------------------------------------------------------------------
# Cardoon (1)
# ---------------------------------------------------------------------
$wsCardoon1 = $wb.Worksheets.Item("Cardoon (1)")

$wsCardoon1.Range("B17").Value = 19.75
$wsCardoon1.Range("C17").Value = 19.5
$wsCardoon1.Range("D17").Value = 6.5
$wsCardoon1.Range("E17").Value = 22.25
$wsCardoon1.Range("F17").Value = 12.5
$wsCardoon1.Range("G17").Value = 14

$wsCardoon1.Range("B18").Value = 20
$wsCardoon1.Range("C18").Value = 19.75
$wsCardoon1.Range("D18").Value = 7
$wsCardoon1.Range("E18").Value = 22.5
$wsCardoon1.Range("F18").Value = 12
$wsCardoon1.Range("G18").Value = 15

$wsCardoon1.Activate()
$wsCardoon1.Range("H23").Select() | Out-Null

# ---------------------------------------------------------------------
# Cardoon (2) - also gains a new "Leaf 7" column (H)
# ---------------------------------------------------------------------
$wsCardoon2 = $wb.Worksheets.Item("Cardoon (2)")

$wsCardoon2.Range("H1").Value = "Leaf 7"

$wsCardoon2.Range("B17").Value = 24.75
$wsCardoon2.Range("D17").Value = 26.5
$wsCardoon2.Range("E17").Value = 24.25
$wsCardoon2.Range("F17").Value = 24.75
$wsCardoon2.Range("G17").Value = 18.25
$wsCardoon2.Range("H17").Value = 2

$wsCardoon2.Range("B18").Value = 24.5
$wsCardoon2.Range("D18").Value = 26.25
$wsCardoon2.Range("E18").Value = 24.25
$wsCardoon2.Range("F18").Value = 25
$wsCardoon2.Range("G18").Value = 20.5
$wsCardoon2.Range("H18").Value = 3.5

$wsCardoon2.Activate()
$wsCardoon2.Range("F22").Select() | Out-Null
